$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Remove the two anchor-table rows (7,8) that no longer exist in the new table ---
$ws.Range('A7:H8').Clear()

# --- Clear old data regions that will be fully rewritten ---
$ws.Range('A3:H6').ClearContents()
$ws.Range('J3:Q34').ClearContents()

# --- Row 1: overall top anchor / type words ---
$ws.Range('A1').Value = 'negative'
$ws.Range('J1').Value = 'positive'

# --- Anchor-score table (A:H), rows 3-6 ---
$ws.Range('A3').Value = 'crude'
$ws.Range('B3').Value = 0.8529411764705882
$ws.Range('C3').Value = 29
$ws.Range('D3').Value = 29
$ws.Range('E3').Value = 0
$ws.Range('F3').Value = 1
$ws.Range('G3').Value = $false
$ws.Range('H3').Value = 5

$ws.Range('A4').Value = 'crisis'
$ws.Range('B4').Value = 0.5993150684931506
$ws.Range('C4').Value = 175
$ws.Range('D4').Value = 175
$ws.Range('E4').Value = 0
$ws.Range('F4').Value = 1
$ws.Range('G4').Value = $false
$ws.Range('H4').Value = 117

$ws.Range('A5').Value = 'panic'
$ws.Range('B5').Value = 0.2189922480620155
$ws.Range('C5').Value = 113
$ws.Range('D5').Value = 113
$ws.Range('E5').Value = 0
$ws.Range('F5').Value = 1
$ws.Range('G5').Value = $false
$ws.Range('H5').Value = 403

$ws.Range('A6').Value = 'sc'
$ws.Range('B6').Value = 0.1693121693121693
$ws.Range('C6').Value = 32
$ws.Range('D6').Value = 32
$ws.Range('E6').Value = 0
$ws.Range('F6').Value = 1
$ws.Range('G6').Value = $false
$ws.Range('H6').Value = 157

# --- Type-occurrences ranking table (J:Q), rows 3-34 ---
$ws.Range('J3').Value = 'love'
$ws.Range('K3').Value = 0.9347826086956522
$ws.Range('L3').Value = 43
$ws.Range('M3').Value = 43
$ws.Range('N3').Value = 1
$ws.Range('O3').Value = 0
$ws.Range('P3').Value = $false
$ws.Range('Q3').Value = 3

$ws.Range('J4').Value = 'best'
$ws.Range('K4').Value = 0.9152542372881356
$ws.Range('L4').Value = 54
$ws.Range('M4').Value = 54
$ws.Range('N4').Value = 1
$ws.Range('O4').Value = 0
$ws.Range('P4').Value = $false
$ws.Range('Q4').Value = 5

$ws.Range('J5').Value = 'interesting'
$ws.Range('K5').Value = 0.8787878787878788
$ws.Range('L5').Value = 29
$ws.Range('M5').Value = 29
$ws.Range('N5').Value = 1
$ws.Range('O5').Value = 0
$ws.Range('P5').Value = $false
$ws.Range('Q5').Value = 4

$ws.Range('J6').Value = 'great'
$ws.Range('K6').Value = 0.8392857142857143
$ws.Range('L6').Value = 94
$ws.Range('M6').Value = 94
$ws.Range('N6').Value = 1
$ws.Range('O6').Value = 0
$ws.Range('P6').Value = $false
$ws.Range('Q6').Value = 18

$ws.Range('J7').Value = 'free'
$ws.Range('K7').Value = 0.8333333333333334
$ws.Range('L7').Value = 100
$ws.Range('M7').Value = 100
$ws.Range('N7').Value = 1
$ws.Range('O7').Value = 0
$ws.Range('P7').Value = $false
$ws.Range('Q7').Value = 20

$ws.Range('J8').Value = 'thanks'
$ws.Range('K8').Value = 0.8292682926829268
$ws.Range('L8').Value = 68
$ws.Range('M8').Value = 68
$ws.Range('N8').Value = 1
$ws.Range('O8').Value = 0
$ws.Range('P8').Value = $false
$ws.Range('Q8').Value = 14

$ws.Range('J9').Value = 'won'
$ws.Range('K9').Value = 0.8205128205128205
$ws.Range('L9').Value = 32
$ws.Range('M9').Value = 32
$ws.Range('N9').Value = 1
$ws.Range('O9').Value = 0
$ws.Range('P9').Value = $false
$ws.Range('Q9').Value = 7

$ws.Range('J10').Value = 'special'
$ws.Range('K10').Value = 0.8055555555555556
$ws.Range('L10').Value = 29
$ws.Range('M10').Value = 29
$ws.Range('N10').Value = 1
$ws.Range('O10').Value = 0
$ws.Range('P10').Value = $false
$ws.Range('Q10').Value = 7

$ws.Range('J11').Value = 'thank'
$ws.Range('K11').Value = 0.8046875
$ws.Range('L11').Value = 103
$ws.Range('M11').Value = 103
$ws.Range('N11').Value = 1
$ws.Range('O11').Value = 0
$ws.Range('P11').Value = $false
$ws.Range('Q11').Value = 25

$ws.Range('J12').Value = 'positive'
$ws.Range('K12').Value = 0.7758620689655172
$ws.Range('L12').Value = 45
$ws.Range('M12').Value = 45
$ws.Range('N12').Value = 1
$ws.Range('O12').Value = 0
$ws.Range('P12').Value = $false
$ws.Range('Q12').Value = 13

$ws.Range('J13').Value = 'confidence'
$ws.Range('K13').Value = 0.75
$ws.Range('L13').Value = 27
$ws.Range('M13').Value = 27
$ws.Range('N13').Value = 1
$ws.Range('O13').Value = 0
$ws.Range('P13').Value = $false
$ws.Range('Q13').Value = 9

$ws.Range('J14').Value = 'safe'
$ws.Range('K14').Value = 0.7394366197183099
$ws.Range('L14').Value = 105
$ws.Range('M14').Value = 105
$ws.Range('N14').Value = 1
$ws.Range('O14').Value = 0
$ws.Range('P14').Value = $false
$ws.Range('Q14').Value = 37

$ws.Range('J15').Value = 'support'
$ws.Range('K15').Value = 0.7075471698113207
$ws.Range('L15').Value = 75
$ws.Range('M15').Value = 75
$ws.Range('N15').Value = 1
$ws.Range('O15').Value = 0
$ws.Range('P15').Value = $false
$ws.Range('Q15').Value = 31

$ws.Range('J16').Value = 'good'
$ws.Range('K16').Value = 0.675
$ws.Range('L16').Value = 108
$ws.Range('M16').Value = 108
$ws.Range('N16').Value = 1
$ws.Range('O16').Value = 0
$ws.Range('P16').Value = $false
$ws.Range('Q16').Value = 52

$ws.Range('J17').Value = 'safety'
$ws.Range('K17').Value = 0.6666666666666666
$ws.Range('L17').Value = 34
$ws.Range('M17').Value = 34
$ws.Range('N17').Value = 1
$ws.Range('O17').Value = 0
$ws.Range('P17').Value = $false
$ws.Range('Q17').Value = 17

$ws.Range('J18').Value = 'better'
$ws.Range('K18').Value = 0.6349206349206349
$ws.Range('L18').Value = 40
$ws.Range('M18').Value = 40
$ws.Range('N18').Value = 1
$ws.Range('O18').Value = 0
$ws.Range('P18').Value = $false
$ws.Range('Q18').Value = 23

$ws.Range('J19').Value = 'fresh'
$ws.Range('K19').Value = 0.625
$ws.Range('L19').Value = 30
$ws.Range('M19').Value = 30
$ws.Range('N19').Value = 1
$ws.Range('O19').Value = 0
$ws.Range('P19').Value = $false
$ws.Range('Q19').Value = 18

$ws.Range('J20').Value = 'well'
$ws.Range('K20').Value = 0.6063829787234043
$ws.Range('L20').Value = 57
$ws.Range('M20').Value = 57
$ws.Range('N20').Value = 1
$ws.Range('O20').Value = 0
$ws.Range('P20').Value = $false
$ws.Range('Q20').Value = 37

$ws.Range('J21').Value = 'relief'
$ws.Range('K21').Value = 0.6
$ws.Range('L21').Value = 30
$ws.Range('M21').Value = 30
$ws.Range('N21').Value = 1
$ws.Range('O21').Value = 0
$ws.Range('P21').Value = $false
$ws.Range('Q21').Value = 20

$ws.Range('J22').Value = 'hand'
$ws.Range('K22').Value = 0.5352480417754569
$ws.Range('L22').Value = 205
$ws.Range('M22').Value = 205
$ws.Range('N22').Value = 1
$ws.Range('O22').Value = 0
$ws.Range('P22').Value = $false
$ws.Range('Q22').Value = 178

$ws.Range('J23').Value = 'care'
$ws.Range('K23').Value = 0.4831460674157304
$ws.Range('L23').Value = 43
$ws.Range('M23').Value = 43
$ws.Range('N23').Value = 1
$ws.Range('O23').Value = 0
$ws.Range('P23').Value = $false
$ws.Range('Q23').Value = 46

$ws.Range('J24').Value = 'like'
$ws.Range('K24').Value = 0.4705882352941176
$ws.Range('L24').Value = 160
$ws.Range('M24').Value = 160
$ws.Range('N24').Value = 1
$ws.Range('O24').Value = 0
$ws.Range('P24').Value = $false
$ws.Range('Q24').Value = 180

$ws.Range('J25').Value = 'help'
$ws.Range('K25').Value = 0.4610169491525424
$ws.Range('L25').Value = 136
$ws.Range('M25').Value = 136
$ws.Range('N25').Value = 1
$ws.Range('O25').Value = 0
$ws.Range('P25').Value = $false
$ws.Range('Q25').Value = 159

$ws.Range('J26').Value = 'protect'
$ws.Range('K26').Value = 0.410958904109589
$ws.Range('L26').Value = 30
$ws.Range('M26').Value = 30
$ws.Range('N26').Value = 1
$ws.Range('O26').Value = 0
$ws.Range('P26').Value = $false
$ws.Range('Q26').Value = 43

$ws.Range('J27').Value = 'increase'
$ws.Range('K27').Value = 0.3846153846153846
$ws.Range('L27').Value = 30
$ws.Range('M27').Value = 30
$ws.Range('N27').Value = 1
$ws.Range('O27').Value = 0
$ws.Range('P27').Value = $false
$ws.Range('Q27').Value = 48

$ws.Range('J28').Value = 'please'
$ws.Range('K28').Value = 0.3514644351464435
$ws.Range('L28').Value = 84
$ws.Range('M28').Value = 84
$ws.Range('N28').Value = 1
$ws.Range('O28').Value = 0
$ws.Range('P28').Value = $false
$ws.Range('Q28').Value = 155

$ws.Range('J29').Value = 'you'
$ws.Range('K29').Value = 0.02416666666666667
$ws.Range('L29').Value = 29
$ws.Range('M29').Value = 29
$ws.Range('N29').Value = 1
$ws.Range('O29').Value = 0
$ws.Range('P29').Value = $false
$ws.Range('Q29').Value = 1171

$ws.Range('J30').Value = '19'
$ws.Range('K30').Value = 0.01215521271622253
$ws.Range('L30').Value = 26
$ws.Range('M30').Value = 30
$ws.Range('N30').Value = 0.87
$ws.Range('O30').Value = 0.13
$ws.Range('P30').Value = $true
$ws.Range('Q30').Value = 2113

$ws.Range('J31').Value = 'and'
$ws.Range('K31').Value = 0.01086549269389284
$ws.Range('L31').Value = 29
$ws.Range('M31').Value = 33
$ws.Range('N31').Value = 0.88
$ws.Range('O31').Value = 0.12
$ws.Range('P31').Value = $true
$ws.Range('Q31').Value = 2640

$ws.Range('J32').Value = 'to'
$ws.Range('K32').Value = 0.007626531083891842
$ws.Range('L32').Value = 33
$ws.Range('M32').Value = 34
$ws.Range('N32').Value = 0.97
$ws.Range('O32').Value = 0.03000000000000003
$ws.Range('P32').Value = $true
$ws.Range('Q32').Value = 4294

$ws.Range('J33').Value = '.'
$ws.Range('K33').Value = 0.007410374524334067
$ws.Range('L33').Value = 37
$ws.Range('M33').Value = 40
$ws.Range('N33').Value = 0.93
$ws.Range('O33').Value = 0.06999999999999995
$ws.Range('P33').Value = $true
$ws.Range('Q33').Value = 4956

$ws.Range('J34').Value = 'the'
$ws.Range('K34').Value = 0.00562015503875969
$ws.Range('L34').Value = 29
$ws.Range('M34').Value = 34
$ws.Range('N34').Value = 0.85
$ws.Range('O34').Value = 0.15
$ws.Range('P34').Value = $true
$ws.Range('Q34').Value = 5131

# --- Re-apply the header/bold+border style ("s=1" in the XML) to new label cells ---
$ws.Range('A2').Copy()
$ws.Range('A3:A6').PasteSpecial(-4122)
$ws.Range('J2').Copy()
$ws.Range('J3:J34').PasteSpecial(-4122)
$excel.CutCopyMode = 0
